$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold font, borders, centered alignment).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I (I0) and J (IF) columns, row by row.
$values = @{
    2  = @(8, 8)
    3  = @(7, 7)
    4  = @(9, 9)
    5  = @(9, 9)
    6  = @(7, 7)
    7  = @(7, 7)
    8  = @(9, 9)
    9  = @(8, 9)
    10 = @(6, 6)
    11 = @(7, 7)
    12 = @(6, 7)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(10, 10)
    16 = @(6, 7)
    17 = @(6, 7)
    18 = @(10, 10)
    19 = @(9, 9)
    20 = @(7, 7)
    21 = @(9, 9)
    22 = @(4, 5)
    23 = @(7, 7)
    24 = @(5, 6)
    25 = @(5, 5)
    26 = @(7, 7)
    27 = @(6, 7)
    28 = @(7, 7)
    29 = @(1, 1)
    30 = @(9, 9)
    31 = @(6, 6)
    32 = @(7, 7)
    33 = @(5, 5)
    34 = @(7, 7)
    35 = @(9, 9)
    36 = @(8, 8)
    37 = @(6, 6)
    38 = @(7, 7)
    39 = @(6, 6)
    40 = @(5, 5)
}

foreach ($row in 2..40) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
